$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B16").Value = "50308ab72f2879be3d0c6f93e42111be"
$ws.Range("B35").Value = "d55e140606222e76082eef082c61f6fd"
$ws.Range("B47").Value = "b326cf6b243488c55c059f20fb6345a4"
$ws.Range("B49").Value = "4826fc67ae504479fe8463172e44d8c8"
$ws.Range("B349").Value = "3d453c9ade000e56d3451cb09588c3b3"
$ws.Range("B360").Value = "06d4b0df2e3daa2e6f3952151324d3c2"
$ws.Range("B361").Value = "9e2b479681aec8331992d5e2e269b068"
$ws.Range("B396").Value = "ca9a0ce7200f67ff0f489c634cd576bf"
$ws.Range("B400").Value = "aec11b26aac47ff6bdcac8864b6f5cbf"
$ws.Range("B405").Value = "363b8da5a43db16b69f56ba299c69d4a"
$ws.Range("B420").Value = "bf3569543f5afe0bd329968445d710df"
$ws.Range("B455").Value = "2d983caf05de9ad5bf2df99f20306b6a"
$ws.Range("B469").Value = "475d7750a415d5eab09d043361d9b844"
$ws.Range("B471").Value = "620b67b1b91269a3195a3efc595edfcf"
$ws.Range("B475").Value = "805427314f487634334aa21bfa53f5ad"
$ws.Range("B645").Value = "bc2165bbd74641c9af2f027e9e9360a9"
$ws.Range("B667").Value = "a42dd747fa3cf21993babb128dfd975c"
$ws.Range("B720").Value = "a35c0b74d2f2f6d675aca19b554f464d"
$ws.Range("B790").Value = "a0eb18d5cc67f9f200f21bc0044efc16"
$ws.Range("B882").Value = "d878f735a89572d2273c1e98708e28dd"
$ws.Range("B961").Value = "1b0fc0e4c5d2c27d6196bfa581be725a"
$ws.Range("B964").Value = "480e3834228ddef3ef8aee5aec97d0f5"
$ws.Range("B974").Value = "67004c0d3bb568ac36e2a173088563ea"
